$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell C1
$ws.Range("C1").Value = "Resolución de imágenes"

# New "Resolución" column values ([144, 144] / "?") for several rows
$ws.Range("C4").Value = "[144, 144]"
$ws.Range("C5").Value = "[144, 144]"
$ws.Range("C7").Value = "[144, 144]"
$ws.Range("C9").Value = "[144, 144]"
$ws.Range("C18").Value = "[144, 144]"
$ws.Range("C16").Value = "?"
$ws.Range("C22").Value = "?"
$ws.Range("C23").Value = "?"

# Rows 22 and 23 changed Interfaz status from "N" to "S"
$ws.Range("B22").Value = "S"
$ws.Range("B23").Value = "S"

# D3 text changed from "P" to "mA"
$ws.Range("D3").Value = "mA"

# New column width for column C
$ws.Columns("C").ColumnWidth = 17

# Update selection / scroll position
$ws.Range("D4").Select()
